$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 652; this shifts the existing rows
# 652-693 down to 654-695 (matching the diff), growing the used range
# from A1:D693 to A1:D695.
$ws.Rows("652:653").Insert()

# Populate the two newly inserted rows with their new data.
# The date-like text in column A must stay literal text (as every other
# date cell in this sheet is stored), so force Text format before
# assigning, then restore the default "Normal" style so no stray
# number-format styling is left on the cell.
$ws.Range("A652").NumberFormat = "@"
$ws.Range("A652").Value = "2026/01/16"
$ws.Range("A652").Style = "Normal"
$ws.Range("B652").Value = "金"
$ws.Range("C652").Value = 22
$ws.Range("D652").Value = 201

$ws.Range("A653").NumberFormat = "@"
$ws.Range("A653").Value = "2026/01/17"
$ws.Range("A653").Style = "Normal"
$ws.Range("B653").Value = "土"
$ws.Range("C653").Value = 2
$ws.Range("D653").Value = 201
